$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 222.36363
$ws.Range("J9").Value = 269.33334
$ws.Range("L9").Value = 269.33334
$ws.Range("N9").Value = -607.33334
$ws.Range("H33").Value = 255
$ws.Range("I33").Value = 249.66667
$ws.Range("J33").Value = 295
$ws.Range("K33").Value = 249.66667
$ws.Range("L33").Value = 295
$ws.Range("M33").Value = -20.66667000000001
$ws.Range("N33").Value = -753
$ws.Range("H98").Value = 770.9583
$ws.Range("I98").Value = 759.2273
$ws.Range("K98").Value = 759.2273
$ws.Range("M98").Value = 738.7727
$ws.Range("H122").Value = 770.9583
$ws.Range("I122").Value = 759.2273
$ws.Range("K122").Value = 2277.6819
$ws.Range("M122").Value = 172.3181
$ws.Range("H139").Value = 42627.145
$ws.Range("J139").Value = 42627.145
$ws.Range("L139").Value = 42627.145
$ws.Range("N139").Value = -52907.145

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 15706.16
$ws.Range("I2").Value = 21870.438
$ws.Range("J2").Value = 4747.4443
$ws.Range("K2").Value = 21870.438
$ws.Range("L2").Value = 4747.4443
$ws.Range("M2").Value = -21757.438
$ws.Range("N2").Value = -4973.4443
$ws.Range("H109").Value = 30785.715
$ws.Range("J109").Value = 30785.715
$ws.Range("L109").Value = 30785.715
$ws.Range("N109").Value = -33559.715
$ws.Range("H116").Value = 15706.16
$ws.Range("I116").Value = 21870.438
$ws.Range("J116").Value = 4747.4443
$ws.Range("K116").Value = 21870.438
$ws.Range("L116").Value = 4747.4443
$ws.Range("M116").Value = -19576.438
$ws.Range("N116").Value = -9335.444299999999
$ws.Range("H122").Value = 5257.9556
$ws.Range("I122").Value = 5168.8423
$ws.Range("K122").Value = 15506.5269
$ws.Range("M122").Value = -13056.5269

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 15706.16
$ws.Range("I3").Value = 21870.438
$ws.Range("J3").Value = 4747.4443
$ws.Range("K3").Value = 21870.438
$ws.Range("L3").Value = 4747.4443
$ws.Range("M3").Value = -21756.438
$ws.Range("N3").Value = -4975.4443
$ws.Range("H54").Value = 2912.4
$ws.Range("I54").Value = 2912.4
$ws.Range("K54").Value = 2912.4
$ws.Range("M54").Value = -2428.4
$ws.Range("H99").Value = 30111.5
$ws.Range("I99").Value = 37524
$ws.Range("J99").Value = 2932.3333
$ws.Range("K99").Value = 37524
$ws.Range("L99").Value = 2932.3333
$ws.Range("M99").Value = -36026
$ws.Range("N99").Value = -5928.3333
$ws.Range("H134").Value = 5304
$ws.Range("I134").Value = 5304
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 15912
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -13377
$ws.Range("N134").ClearContents()

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3071.577
$ws.Range("I31").Value = 2414
$ws.Range("K31").Value = 2414
$ws.Range("M31").Value = -2119
$ws.Range("H34").Value = 3071.577
$ws.Range("I34").Value = 2414
$ws.Range("K34").Value = 2414
$ws.Range("M34").Value = -2212
$ws.Range("H74").Value = 69656.5
$ws.Range("J74").Value = 69656.5
$ws.Range("L74").Value = 69656.5
$ws.Range("N74").Value = -71404.5
$ws.Range("H77").Value = 69656.5
$ws.Range("J77").Value = 69656.5
$ws.Range("L77").Value = 208969.5
$ws.Range("N77").Value = -217705.5
$ws.Range("H109").Value = 12888.667
$ws.Range("J109").Value = 12888.667
$ws.Range("L109").Value = 12888.667
$ws.Range("N109").Value = -14968.667
$ws.Range("H134").Value = 4859.759
$ws.Range("I134").Value = 4823.846
$ws.Range("J134").Value = 5171
$ws.Range("K134").Value = 14471.538
$ws.Range("L134").Value = 15513
$ws.Range("M134").Value = -11936.538
$ws.Range("N134").Value = -20583

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H126").Value = 2499
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 2499
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 7497
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -17377

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3034375.5
$ws.Range("I97").Value = 582.9048
$ws.Range("J97").Value = 10113225
$ws.Range("K97").Value = 582.9048
$ws.Range("L97").Value = 10113225
$ws.Range("M97").Value = -86.90480000000002
$ws.Range("N97").Value = -10114217
$ws.Range("H108").Value = 100001
$ws.Range("J108").Value = 100001
$ws.Range("L108").Value = 100001
$ws.Range("N108").Value = -107681
$ws.Range("H113").Value = 5900.8125
$ws.Range("I113").Value = 4311.143
$ws.Range("K113").Value = 4311.143
$ws.Range("M113").Value = -2141.143
$ws.Range("H122").Value = 4435.857
$ws.Range("J122").Value = 4941.727
$ws.Range("L122").Value = 14825.181
$ws.Range("N122").Value = -19725.181
$ws.Range("H126").Value = 7786.9473
$ws.Range("I126").Value = 6175.6665
$ws.Range("J126").Value = 8530.615
$ws.Range("K126").Value = 18526.9995
$ws.Range("L126").Value = 25591.845
$ws.Range("M126").Value = -16056.9995
$ws.Range("N126").Value = -30531.845

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3125.8438
$ws.Range("I7").Value = 3139.7307
$ws.Range("J7").Value = 3065.6667
$ws.Range("K7").Value = 3139.7307
$ws.Range("L7").Value = 3065.6667
$ws.Range("M7").Value = -3027.7307
$ws.Range("N7").Value = -3289.6667
$ws.Range("H22").Value = 1156.25
$ws.Range("J22").Value = 875
$ws.Range("L22").Value = 875
$ws.Range("N22").Value = -1465
$ws.Range("H27").Value = 1156.25
$ws.Range("J27").Value = 875
$ws.Range("L27").Value = 875
$ws.Range("N27").Value = -1089
$ws.Range("H61").Value = 1996.3334
$ws.Range("I61").Value = 1997
$ws.Range("J61").Value = 1995
$ws.Range("K61").Value = 1997
$ws.Range("L61").Value = 1995
$ws.Range("M61").Value = -1795
$ws.Range("N61").Value = -2399
$ws.Range("H100").Value = 194051.19
$ws.Range("I100").Value = 194051.19
$ws.Range("K100").Value = 194051.19
$ws.Range("M100").Value = -193510.19
$ws.Range("H112").Value = 32285.643
$ws.Range("J112").Value = 32285.643
$ws.Range("L112").Value = 32285.643
$ws.Range("N112").Value = -35239.643
$ws.Range("H113").Value = 1996.3334
$ws.Range("I113").Value = 1997
$ws.Range("J113").Value = 1995
$ws.Range("K113").Value = 1997
$ws.Range("L113").Value = 1995
$ws.Range("M113").Value = 173
$ws.Range("N113").Value = -6335
$ws.Range("H122").Value = 4585.1875
$ws.Range("I122").Value = 3696.6365
$ws.Range("J122").Value = 6540
$ws.Range("K122").Value = 11089.9095
$ws.Range("L122").Value = 19620
$ws.Range("M122").Value = -8639.9095
$ws.Range("N122").Value = -24520
$ws.Range("H126").Value = 3125.8438
$ws.Range("I126").Value = 3139.7307
$ws.Range("J126").Value = 3065.6667
$ws.Range("K126").Value = 9419.1921
$ws.Range("L126").Value = 9197.000100000001
$ws.Range("M126").Value = -6949.1921
$ws.Range("N126").Value = -14137.0001

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 59248.5
$ws.Range("J27").Value = 59248.5
$ws.Range("L27").Value = 59248.5
$ws.Range("N27").Value = -59386.5
$ws.Range("H101").Value = 42500
$ws.Range("J101").Value = 42500
$ws.Range("L101").Value = 42500
$ws.Range("N101").Value = -48990
$ws.Range("H107").Value = 1038.579
$ws.Range("I107").Value = 1193
$ws.Range("J107").Value = 773.8570999999999
$ws.Range("K107").Value = 3579
$ws.Range("L107").Value = 2321.5713
$ws.Range("M107").Value = -1659
$ws.Range("N107").Value = -6161.5713
$ws.Range("H109").Value = 13722.223
$ws.Range("J109").Value = 13722.223
$ws.Range("L109").Value = 13722.223
$ws.Range("N109").Value = -16496.223
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H115").Value = 21999.916
$ws.Range("J115").Value = 21999.916
$ws.Range("L115").Value = 21999.916
$ws.Range("N115").Value = -25133.916
$ws.Range("H132").Value = 4305.2
$ws.Range("I132").Value = 4313.6924
$ws.Range("J132").Value = 4250
$ws.Range("K132").Value = 12941.0772
$ws.Range("L132").Value = 12750
$ws.Range("M132").Value = -10411.0772
$ws.Range("N132").Value = -17810
